# Fruta / hortaliza, semanal
# Insert a new weekly record at row 82 (pushing the existing rows 82-94
# down to 83-95) and populate it with the new sample.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(82).Insert()

$ws.Cells.Item(82, 1).Value  = 5
$ws.Cells.Item(82, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(82, 3).Value  = "Maule"
$ws.Cells.Item(82, 4).Value  = 44841
$ws.Cells.Item(82, 5).Value  = 7
$ws.Cells.Item(82, 6).Value  = 100112026
$ws.Cells.Item(82, 7).Value  = "Haba"
$ws.Cells.Item(82, 8).Value  = "Sin especificar"
$ws.Cells.Item(82, 9).Value  = "Primera"
$ws.Cells.Item(82, 10).Value = 200
$ws.Cells.Item(82, 11).Value = 9000
$ws.Cells.Item(82, 12).Value = 9000
$ws.Cells.Item(82, 13).Value = 9000
$ws.Cells.Item(82, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(82, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(82, 16).Value = 360
$ws.Cells.Item(82, 17).Value = 25
$ws.Cells.Item(82, 18).Value = "Hortaliza"
